# "tried more anti-overfitting methods in random forest"
#
# - renames the second sheet "Feuil3" -> "results_tfidf"
# - populates that sheet with a results_tfidf table mirroring the layout
#   already used on results_allgenres (Logistic / SVM blocks), adding a
#   third x_train variant ("x_train_pca (10 components, 28% explained var)")
# - updates the selection on both sheets

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # results_allgenres
$ws2 = $wb.Worksheets.Item(2)   # Feuil3 -> results_tfidf

# --- rename the second sheet -------------------------------------------------
$ws2.Name = "results_tfidf"

# --- bring over the formatted header blocks ----------------------------------
# results_allgenres already has identically-styled "Logistic" (B1:E3) and
# "SVM" (J1:M3) header blocks (merged cells, fills, centered headers); reuse
# them verbatim instead of re-building the formatting from scratch.
$ws1.Range("B1:E3").Copy($ws2.Range("B1:E3"))
$ws1.Range("J1:M3").Copy($ws2.Range("F1:I3"))

# --- row labels ---------------------------------------------------------------
$ws2.Range("A4").Value = "x_train"
$ws2.Range("A5").Value = "x_train2(0-1k)"
$ws2.Range("A6").Value = "x_train_pca (10 components, 28% explained var)"

# --- data: precision_1/recall_1 x test/train x Logistic/SVM ------------------
$data = @(
  @(0.07, 0.42, 0.1,  0.56, 0.07, 0.39, 0.12, 0.62),
  @(0.07, 0.43, 0.09, 0.55, 0.07, 0.42, 0.1,  0.59),
  @(0.07, 0.42, 0.07, 0.44, 0.07, 0.32, 0.07, 0.33)
)
for ($r = 0; $r -lt $data.Length; $r++) {
  $row = 4 + $r
  for ($c = 0; $c -lt 8; $c++) {
    $col = 2 + $c
    $ws2.Cells.Item($row, $col).Value = $data[$r][$c]
  }
}

# --- column A width, matching results_allgenres's label column ---------------
$ws2.Columns.Item(1).ColumnWidth = $ws1.Columns.Item(1).ColumnWidth

# --- selections: set sheet2's remembered selection first, then re-select on
#     sheet1 last so results_allgenres stays the active/visible tab ----------
$ws2.Range("C7").Select()
$ws1.Range("A5").Select()
